$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.134.59"
$ws.Range("E2").Value = "  -3.26%  "
$ws.Range("D3").Value = "'1.860.53"
$ws.Range("E3").Value = "  -4.39%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'233.74"
$ws.Range("E5").Value = "  -3.69%  "
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "'0.4649"
$ws.Range("E7").Value = "  -3.46%  "
$ws.Range("D8").Value = "'0.2818"
$ws.Range("E8").Value = "  -3.37%  "
$ws.Range("D9").Value = "'0.06546"
$ws.Range("E9").Value = "  -3.84%  "
$ws.Range("D10").Value = "'19.83"
$ws.Range("E10").Value = "  -1.40%  "
$ws.Range("D11").Value = "'0.07809"
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").Value = "'96.42"
$ws.Range("E12").Value = "  -7.79%  "
$ws.Range("D13").Value = "'1.857.98"
$ws.Range("E13").Value = "  -4.61%  "
$ws.Range("D14").Value = "'5.126"
$ws.Range("E14").Value = "  -3.44%  "
$ws.Range("D15").Value = "'0.6663"
$ws.Range("E15").Value = "  -3.42%  "
$ws.Range("D16").Value = "'281.06"
$ws.Range("E16").Value = "  -5.44%  "
$ws.Range("D17").Value = "'30.167.53"
$ws.Range("E17").Value = "  -3.19%  "
$ws.Range("D18").Value = "'0.9992"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").Value = "'5.516"
$ws.Range("E19").Value = "  -1.37%  "
$ws.Range("D20").Value = "'12.61"
$ws.Range("E20").Value = "  -3.12%  "
$ws.Range("D21").Value = "'2.101.18"
$ws.Range("E21").Value = "  -5.02%  "
$ws.Range("D22").Value = "'0.000007239"
$ws.Range("E22").Value = "  -5.00%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "'6.126"
$ws.Range("E24").Value = "  -5.26%  "
$ws.Range("D25").Value = "'9.318"
$ws.Range("E25").Value = "  -2.81%  "
$ws.Range("D26").Value = "'166.17"
$ws.Range("E26").Value = "  -1.68%  "
$ws.Range("D27").Value = "'18.88"
$ws.Range("E27").Value = "  -5.03%  "
$ws.Range("D28").Value = "'1.909"
$ws.Range("E28").Value = "  -10.71%  "
$ws.Range("D29").Value = "'1.339"
$ws.Range("E29").Value = "  -3.71%  "
$ws.Range("D30").Value = "'0.09574"
$ws.Range("E30").Value = "  -5.76%  "
$ws.Range("D31").Value = "'4.402"
$ws.Range("E31").Value = "  -5.07%  "
$ws.Range("D32").Value = "'1.469"
$ws.Range("E32").Value = "  -4.58%  "
$ws.Range("D33").Value = "'4.110"
$ws.Range("E33").Value = "  -5.75%  "
$ws.Range("D34").Value = "'0.04657"
$ws.Range("E34").Value = "  -3.81%  "
$ws.Range("D35").Value = "'0.7003"
$ws.Range("E35").Value = "  -5.77%  "
$ws.Range("D36").Value = "'1.092"
$ws.Range("E36").Value = "  -3.73%  "
$ws.Range("D37").Value = "'2.702"
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("D38").Value = "'0.01851"
$ws.Range("E38").Value = "  -5.54%  "
$ws.Range("D39").Value = "'6.276"
$ws.Range("E39").Value = "  -4.94%  "
$ws.Range("D40").Value = "'2.518"
$ws.Range("E40").Value = "  -4.94%  "
$ws.Range("D41").Value = "'73.36"
$ws.Range("E41").Value = "  -5.80%  "
$ws.Range("D42").Value = "'0.8536"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("D43").Value = "'1.918"
$ws.Range("E43").Value = "  -5.82%  "
$ws.Range("D44").Value = "'0.9996"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "'0.4154"
$ws.Range("E45").Value = "  -5.13%  "
$ws.Range("D46").Value = "'103.46"
$ws.Range("E46").Value = "  -2.89%  "
$ws.Range("D47").Value = "'990.17"
$ws.Range("E47").Value = "  -3.70%  "
$ws.Range("D48").Value = "'7.190"
$ws.Range("E48").Value = "  -5.18%  "
$ws.Range("D49").Value = "'9.232"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").Value = "'34.17"
$ws.Range("E50").Value = "  -3.05%  "
$ws.Range("D51").Value = "'0.1139"
$ws.Range("E51").Value = "  -6.25%  "
